$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 6
$ws.Range("H6").Value = 48.333332
$ws.Range("I6").Value = 48.333332
$ws.Range("K6").Value = 144.999996
$ws.Range("M6").Value = -32.99999600000001
# row 17
$ws.Range("H17").Value = 4255736
$ws.Range("J17").Value = 4440685
$ws.Range("L17").Value = 13322055
$ws.Range("N17").Value = -13322391
# row 51
$ws.Range("H51").Value = 6250
$ws.Range("J51").Value = 5000
$ws.Range("L51").Value = 5000
$ws.Range("N51").Value = -5968
# row 132
$ws.Range("H132").Value = 4017.261
$ws.Range("I132").Value = 4999.1177
$ws.Range("J132").Value = 1235.3334
$ws.Range("K132").Value = 14997.3531
$ws.Range("L132").Value = 3706.0002
$ws.Range("M132").Value = -12467.3531
$ws.Range("N132").Value = -8766.0002
# row 137
$ws.Range("H137").Value = 2018.75
$ws.Range("I137").Value = 2180
$ws.Range("J137").Value = 1750
$ws.Range("K137").Value = 6540
$ws.Range("L137").Value = 5250
$ws.Range("M137").Value = -3990
$ws.Range("N137").Value = -10350
# row 138
$ws.Range("H138").Value = 27029820
$ws.Range("J138").Value = 3456.5454
$ws.Range("L138").Value = 10369.6362
$ws.Range("N138").Value = -20649.6362

$ws = $wb.Worksheets.Item("ARM")
# row 63
$ws.Range("H63").Value = 2049.7693
$ws.Range("I63").Value = 2068.0908
$ws.Range("J63").Value = 1949
$ws.Range("K63").Value = 2068.0908
$ws.Range("L63").Value = 1949
$ws.Range("M63").Value = -1382.0908
$ws.Range("N63").Value = -3321
# row 66
$ws.Range("H66").Value = 2049.7693
$ws.Range("I66").Value = 2068.0908
$ws.Range("J66").Value = 1949
$ws.Range("K66").Value = 10340.454
$ws.Range("L66").Value = 9745
$ws.Range("M66").Value = -6908.454
$ws.Range("N66").Value = -16609
# row 74
$ws.Range("H74").Value = 47620932
$ws.Range("I74").Value = 111111740
$ws.Range("J74").Value = 2822
$ws.Range("K74").Value = 111111740
$ws.Range("L74").Value = 2822
$ws.Range("M74").Value = -111110866
$ws.Range("N74").Value = -4570
# row 77
$ws.Range("H77").Value = 47620932
$ws.Range("I77").Value = 111111740
$ws.Range("J77").Value = 2822
$ws.Range("K77").Value = 555558700
$ws.Range("L77").Value = 14110
$ws.Range("M77").Value = -555554332
$ws.Range("N77").Value = -22846
# row 133
$ws.Range("H133").Value = 180000
$ws.Range("J133").Value = 180000
$ws.Range("L133").Value = 180000
$ws.Range("N133").Value = -185060

$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 3312.9534
$ws.Range("I31").Value = 2882.2856
$ws.Range("K31").Value = 2882.2856
$ws.Range("M31").Value = -2587.2856
# row 34
$ws.Range("H34").Value = 3312.9534
$ws.Range("I34").Value = 2882.2856
$ws.Range("K34").Value = 2882.2856
$ws.Range("M34").Value = -2680.2856
# row 58
$ws.Range("H58").Value = 20991.154
$ws.Range("I58").Value = 1885.8182
$ws.Range("K58").Value = 1885.8182
$ws.Range("M58").Value = -1682.8182
# row 136
$ws.Range("H136").Value = 20991.154
$ws.Range("I136").Value = 1885.8182
$ws.Range("K136").Value = 5657.4546
$ws.Range("M136").Value = -3107.4546

$ws = $wb.Worksheets.Item("CUL")
# row 5
$ws.Range("H5").Value = 1801.4
$ws.Range("I5").Value = 1130.8
$ws.Range("J5").Value = 2472
$ws.Range("K5").Value = 3392.4
$ws.Range("L5").Value = 7416
$ws.Range("M5").Value = -3280.4
$ws.Range("N5").Value = -7640
# row 107
$ws.Range("H107").Value = 4937.619
$ws.Range("I107").Value = 5726.6665
$ws.Range("J107").Value = 203.33333
$ws.Range("K107").Value = 17179.9995
$ws.Range("L107").Value = 609.99999
$ws.Range("M107").Value = -15259.9995
$ws.Range("N107").Value = -4449.99999
# row 113
$ws.Range("H113").Value = 719.13336
$ws.Range("I113").Value = 626.3333
$ws.Range("J113").Value = 781
$ws.Range("K113").Value = 1878.9999
$ws.Range("L113").Value = 2343
$ws.Range("M113").Value = 291.0001
$ws.Range("N113").Value = -6683
# row 131
$ws.Range("H131").Value = 744.5714
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 744.5714
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 2233.7142
$ws.Range("M131").Value = $null
$ws.Range("N131").Value = -12313.7142
# row 135
$ws.Range("H135").Value = 1801.4
$ws.Range("I135").Value = 1130.8
$ws.Range("J135").Value = 2472
$ws.Range("K135").Value = 10177.2
$ws.Range("L135").Value = 22248
$ws.Range("M135").Value = -7642.199999999999
$ws.Range("N135").Value = -27318
# row 136
$ws.Range("H136").Value = 3298.6
$ws.Range("I136").Value = 997.5
$ws.Range("J136").Value = 4832.6665
$ws.Range("K136").Value = 2992.5
$ws.Range("L136").Value = 14497.9995
$ws.Range("M136").Value = 2107.5
$ws.Range("N136").Value = -24697.9995

$ws = $wb.Worksheets.Item("GSM")
# row 102
$ws.Range("H102").Value = 16669001
$ws.Range("I102").Value = 19232914
$ws.Range("J102").Value = 3565.75
$ws.Range("K102").Value = 19232914
$ws.Range("L102").Value = 3565.75
$ws.Range("M102").Value = -19231292
$ws.Range("N102").Value = -6809.75
# row 107
$ws.Range("H107").Value = 2849197.2
$ws.Range("I107").Value = 183.41176
$ws.Range("J107").Value = 7692520.5
$ws.Range("K107").Value = 183.41176
$ws.Range("L107").Value = 7692520.5
$ws.Range("M107").Value = 1736.58824
$ws.Range("N107").Value = -7696360.5
# row 113
$ws.Range("H113").Value = 6984.4707
$ws.Range("I113").Value = 9335.091
$ws.Range("J113").Value = 2675
$ws.Range("K113").Value = 9335.091
$ws.Range("L113").Value = 2675
$ws.Range("M113").Value = -7165.091
$ws.Range("N113").Value = -7015
# row 123
$ws.Range("H123").Value = 5066.625
$ws.Range("I123").Value = 3267.1428
$ws.Range("J123").Value = 17663
$ws.Range("K123").Value = 3267.1428
$ws.Range("L123").Value = 17663
$ws.Range("M123").Value = -817.1428000000001
$ws.Range("N123").Value = -22563

$ws = $wb.Worksheets.Item("LTW")
# row 82
$ws.Range("H82").Value = 2150
$ws.Range("I82").Value = 2340
$ws.Range("J82").Value = 1200
$ws.Range("K82").Value = 2340
$ws.Range("L82").Value = 1200
$ws.Range("M82").Value = -1979
$ws.Range("N82").Value = -1922
# row 85
$ws.Range("H85").Value = 2150
$ws.Range("I85").Value = 2340
$ws.Range("J85").Value = 1200
$ws.Range("K85").Value = 2340
$ws.Range("L85").Value = 1200
$ws.Range("M85").Value = -1092
$ws.Range("N85").Value = -3696
# row 136
$ws.Range("H136").Value = 1278.7894
$ws.Range("I136").Value = 1299.8334
$ws.Range("J136").Value = 900
$ws.Range("K136").Value = 3899.5002
$ws.Range("L136").Value = 2700
$ws.Range("M136").Value = -1349.5002
$ws.Range("N136").Value = -7800

$ws = $wb.Worksheets.Item("WVR")
# row 3
$ws.Range("H3").Value = 400
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 400
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 400
$ws.Range("M3").Value = $null
$ws.Range("N3").Value = -628
